$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": add column AF for the new wave "22. 6. 2021"
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Copy the header style from the previous date column (AE1) so the new
# header cell (AF1) keeps the same bold/border/alignment formatting.
$wsData.Range("AE1").Copy($wsData.Range("AF1"))
$wsData.Range("AF1").Value = "22. 6. 2021"

$dataValues = @{
    2 = 0.59
    3 = 0.25
    4 = 0.16
    5 = 0.75
    6 = 0.16
    7 = 0.09
    8 = 0.64
    9 = 0.23
    10 = 0.13
    11 = 0.43
    12 = 0.33
    13 = 0.24
    14 = 0.38
    15 = 0.3
    16 = 0.32
    17 = 0.5600000000000001
    18 = 0.28
    19 = 0.16
    20 = 0.65
    21 = 0.23
    22 = 0.12
    23 = 0.75
    24 = 0.17
    25 = 0.08
    26 = 0.61
    27 = 0.25
    28 = 0.14
    29 = 0.38
    30 = 0.35
    31 = 0.27
    32 = 0.6
    33 = 0.29
    34 = 0.11
    35 = 0.55
    36 = 0.25
    37 = 0.2
    38 = 0.6
    39 = 0.24
    40 = 0.16
    41 = 0.6
    42 = 0.23
    43 = 0.17
    44 = 0.58
    45 = 0.27
    46 = 0.15
    47 = 0.62
    48 = 0.21
    49 = 0.17
    50 = 0.55
    51 = 0.3
    52 = 0.15
    53 = 0.5600000000000001
    54 = 0.26
    55 = 0.18
    56 = 0.54
    57 = 0.33
    58 = 0.13
}

foreach ($row in $dataValues.Keys) {
    $wsData.Cells.Item($row, 32).Value = $dataValues[$row]
}

# Update the title strings (now referencing "28. 6. 2021" instead of "1. 6. 2021")
$wsData.Range("A59").Value = "Život během pandemie, Počet protektivních aktivit, % respondentů celkově a ve skupinách, aktualizace 28. 6. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": add column AE for the new wave "22. 6. 2021"
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

$wsPocet.Range("AD1").Copy($wsPocet.Range("AE1"))
$wsPocet.Range("AE1").Value = "22. 6. 2021"

$pocetValues = @{
    2 = 1904
    3 = 459
    4 = 699
    5 = 746
    6 = 304
    7 = 541
    8 = 844
    9 = 573
    10 = 803
    11 = 528
    12 = 377
    13 = 420
    14 = 1107
    15 = 937
    16 = 967
    17 = 983
    18 = 447
    19 = 223
    20 = 251
}

foreach ($row in $pocetValues.Keys) {
    $wsPocet.Cells.Item($row, 31).Value = $pocetValues[$row]
}

# Trailing blank placeholder cell, matching the rest of row 21
$wsPocet.Range("AD21").Copy($wsPocet.Range("AE21"))

# Update the title strings (now referencing "28. 6. 2021" instead of "1. 6. 2021")
$wsPocet.Range("A21").Value = "Život během pandemie, Počet protektivních aktivit, velikost dotázaného souboru celkově a ve skupinách, aktualizace 28. 6. 2021"
